# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 138
    3  = 1665
    4  = 646
    5  = 1112
    7  = 11696
    10 = 465
    11 = 384
    12 = 1099
    13 = 821
    14 = 13402
    15 = 13241
    17 = 147
    20 = 259
    23 = 140
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
